# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# Updates columns H-N (price/profit figures) for specific rows across all 8 sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 2500.389
$ws.Range("I118").Value = 1510
$ws.Range("J118").Value = 3738.375
$ws.Range("K118").Value = 4530
$ws.Range("L118").Value = 11215.125
$ws.Range("M118").Value = -2873
$ws.Range("N118").Value = -14529.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6792.3384
$ws.Range("I32").Value = 6947.0557
$ws.Range("J32").Value = 6032.8184
$ws.Range("K32").Value = 6947.0557
$ws.Range("L32").Value = 6032.8184
$ws.Range("M32").Value = -6660.0557
$ws.Range("N32").Value = -6606.8184

$ws.Range("H70").Value = 37859
$ws.Range("J70").Value = 37859
$ws.Range("L70").Value = 37859
$ws.Range("N70").Value = -38399

$ws.Range("H73").Value = 37859
$ws.Range("J73").Value = 37859
$ws.Range("L73").Value = 37859
$ws.Range("N73").Value = -39731

$ws.Range("H74").Value = 2281.6296
$ws.Range("I74").Value = 1356.7059
$ws.Range("J74").Value = 3854
$ws.Range("K74").Value = 1356.7059
$ws.Range("L74").Value = 3854
$ws.Range("M74").Value = -482.7058999999999
$ws.Range("N74").Value = -5602

$ws.Range("H77").Value = 2281.6296
$ws.Range("I77").Value = 1356.7059
$ws.Range("J77").Value = 3854
$ws.Range("K77").Value = 6783.5295
$ws.Range("L77").Value = 19270
$ws.Range("M77").Value = -2415.5295
$ws.Range("N77").Value = -28006

$ws.Range("H97").Value = 4679.933
$ws.Range("I97").Value = 5144.36
$ws.Range("J97").Value = 2357.8
$ws.Range("K97").Value = 5144.36
$ws.Range("L97").Value = 2357.8
$ws.Range("M97").Value = -4648.36
$ws.Range("N97").Value = -3349.8

$ws.Range("H102").Value = 5327.16
$ws.Range("I102").Value = 5507.7827
$ws.Range("J102").Value = 3250
$ws.Range("K102").Value = 5507.7827
$ws.Range("L102").Value = 3250
$ws.Range("M102").Value = -3885.7827
$ws.Range("N102").Value = -6494

$ws.Range("H135").Value = 25196.5
$ws.Range("J135").Value = 25196.5
$ws.Range("L135").Value = 25196.5
$ws.Range("N135").Value = -35336.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J22").Value = 200
$ws.Range("L22").Value = 200
$ws.Range("N22").Value = -546

$ws.Range("H64").Value = 795.0769
$ws.Range("I64").Value = 1198.1428
$ws.Range("J64").Value = 324.83334
$ws.Range("K64").Value = 1198.1428
$ws.Range("L64").Value = 324.83334
$ws.Range("M64").Value = -973.1428000000001
$ws.Range("N64").Value = -774.83334

$ws.Range("H67").Value = 795.0769
$ws.Range("I67").Value = 1198.1428
$ws.Range("J67").Value = 324.83334
$ws.Range("K67").Value = 1198.1428
$ws.Range("L67").Value = 324.83334
$ws.Range("M67").Value = -418.1428000000001
$ws.Range("N67").Value = -1884.83334

$ws.Range("H94").Value = 1301.3334
$ws.Range("I94").Value = 717.1111
$ws.Range("J94").Value = 1885.5555
$ws.Range("K94").Value = 717.1111
$ws.Range("L94").Value = 1885.5555
$ws.Range("M94").Value = -266.1111
$ws.Range("N94").Value = -2787.5555

$ws.Range("H99").Value = 1621.4762
$ws.Range("I99").Value = 1532.2142
$ws.Range("J99").Value = 1800
$ws.Range("K99").Value = 1532.2142
$ws.Range("L99").Value = 1800
$ws.Range("M99").Value = -34.21419999999989
$ws.Range("N99").Value = -4796

$ws.Range("H107").Value = 2094.524
$ws.Range("I107").Value = 2370.8572
$ws.Range("K107").Value = 2370.8572
$ws.Range("M107").Value = -450.8571999999999

$ws.Range("H134").Value = 2987.8948
$ws.Range("I134").Value = 1978.6086
$ws.Range("J134").Value = 7208.5454
$ws.Range("K134").Value = 5935.825800000001
$ws.Range("L134").Value = 21625.6362
$ws.Range("M134").Value = -3400.825800000001
$ws.Range("N134").Value = -26695.6362

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1000
$ws.Range("I16").Value = 1000
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 1000
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -713
$ws.Range("N16").Value = -1574

$ws.Range("H113").Value = 1000
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 1170
$ws.Range("N113").Value = -5340

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 1010140.2
$ws.Range("I33").Value = 1122374.5
$ws.Range("J33").Value = 32
$ws.Range("K33").Value = 6734247
$ws.Range("L33").Value = 192
$ws.Range("M33").Value = -6733964
$ws.Range("N33").Value = -758

$ws.Range("H63").Value = 1530
$ws.Range("J63").Value = 7000
$ws.Range("L63").Value = 21000
$ws.Range("N63").Value = -22498

$ws.Range("H66").Value = 1530
$ws.Range("J66").Value = 7000
$ws.Range("L66").Value = 63000
$ws.Range("N66").Value = -70488

$ws.Range("H80").Value = 2816.5833
$ws.Range("J80").Value = 3757
$ws.Range("L80").Value = 11271
$ws.Range("N80").Value = -13143

$ws.Range("H81").Value = 1086.1428
$ws.Range("I81").Value = 520.6
$ws.Range("J81").Value = 2500
$ws.Range("K81").Value = 1561.8
$ws.Range("L81").Value = 7500
$ws.Range("M81").Value = -438.8000000000002
$ws.Range("N81").Value = -9746

$ws.Range("H83").Value = 2816.5833
$ws.Range("J83").Value = 3757
$ws.Range("L83").Value = 33813
$ws.Range("N83").Value = -43173

$ws.Range("H84").Value = 1086.1428
$ws.Range("I84").Value = 520.6
$ws.Range("J84").Value = 2500
$ws.Range("K84").Value = 4685.400000000001
$ws.Range("L84").Value = 22500
$ws.Range("M84").Value = 930.5999999999995
$ws.Range("N84").Value = -33732

$ws.Range("H87").Value = 17053.334
$ws.Range("I87").Value = 5300
$ws.Range("J87").Value = 19991.666
$ws.Range("K87").Value = 15900
$ws.Range("L87").Value = 59974.99800000001
$ws.Range("M87").Value = -14652
$ws.Range("N87").Value = -62470.99800000001

$ws.Range("H90").Value = 17053.334
$ws.Range("I90").Value = 5300
$ws.Range("J90").Value = 19991.666
$ws.Range("K90").Value = 47700
$ws.Range("L90").Value = 179924.994
$ws.Range("M90").Value = -41460
$ws.Range("N90").Value = -192404.994

$ws.Range("H92").Value = 486.1111
$ws.Range("I92").Value = 295.5
$ws.Range("J92").Value = 638.6
$ws.Range("K92").Value = 886.5
$ws.Range("L92").Value = 1915.8
$ws.Range("M92").Value = 361.5
$ws.Range("N92").Value = -4411.8

$ws.Range("H107").Value = 861.53424
$ws.Range("I107").Value = 319
$ws.Range("J107").Value = 1161.6595
$ws.Range("K107").Value = 957
$ws.Range("L107").Value = 3484.9785
$ws.Range("M107").Value = 963
$ws.Range("N107").Value = -7324.9785

$ws.Range("H108").Value = 1940.8334
$ws.Range("I108").Value = 829
$ws.Range("J108").Value = 7500
$ws.Range("K108").Value = 2487
$ws.Range("L108").Value = 22500
$ws.Range("M108").Value = 393
$ws.Range("N108").Value = -28260

$ws.Range("H109").Value = 3600
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 3600
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 10800
$ws.Range("N109").Value = -12880
$ws.Range("M109").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 533.3333
$ws.Range("J13").Value = 667
$ws.Range("L13").Value = 667
$ws.Range("N13").Value = -945

$ws.Range("H113").Value = 1495.0834
$ws.Range("I113").Value = 1512.8182
$ws.Range("K113").Value = 1512.8182
$ws.Range("M113").Value = 657.1818000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 583.4
$ws.Range("I55").Value = 299.75
$ws.Range("J55").Value = 686.5454999999999
$ws.Range("K55").Value = 299.75
$ws.Range("L55").Value = 686.5454999999999
$ws.Range("M55").Value = -126.75
$ws.Range("N55").Value = -1032.5455

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3190.4424
$ws.Range("I132").Value = 3328
$ws.Range("J132").Value = 2306.1428
$ws.Range("K132").Value = 9984
$ws.Range("L132").Value = 6918.428400000001
$ws.Range("M132").Value = -7454
$ws.Range("N132").Value = -11978.4284

$ws.Range("H136").Value = 3704687.5
$ws.Range("I136").Value = 4000902.2
$ws.Range("J136").Value = 2002.5
$ws.Range("K136").Value = 12002706.6
$ws.Range("L136").Value = 6007.5
$ws.Range("M136").Value = -12000156.6
$ws.Range("N136").Value = -11107.5

$ws.Range("H137").Value = 58709.223
$ws.Range("J137").Value = 58709.223
$ws.Range("L137").Value = 58709.223
$ws.Range("N137").Value = -68909.223

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
